# Regenerate the handback-status report values (re-run of the handback
# pipeline): the source doc "2f49caf1-...md" was replaced by a new GUID-named
# doc "dae40f3e-...md" and "e1908cf4-...md" was replaced by "ffff6a949d48-...md",
# a fresh handoff/handback cycle ran, and the Correspond Handoff/Handback xlf
# file names + timestamps were updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$ov.Range("B2").Value = "e2e\dae40f3e-36d2-4979-815f-7582a5f37575.md"
$ov.Range("G2").Value = "2016-09-06 17:44:27"

$ov.Range("A3").Value = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$ov.Range("B3").Value = "e2e\ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$ov.Range("G3").Value = "2016-09-06 17:44:27"

# keep the hyperlink captions in sync with the new file names
$ov.Hyperlinks.Item(1).TextToDisplay = "e2e\dae40f3e-36d2-4979-815f-7582a5f37575.md"
$ov.Hyperlinks.Item(2).TextToDisplay = "e2e\ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"

# ---------------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$zh.Range("G2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.zh-cn.xlf"
$zh.Range("H2").Value = "2016-09-06 17:44:22"
$zh.Range("I2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$zh.Range("J2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-06 17:44:40"

$zh.Range("A3").Value = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$zh.Range("G3").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-06 17:44:22"
$zh.Range("I3").Value = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$zh.Range("J3").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-06 17:44:40"

$zh.Hyperlinks.Item(1).TextToDisplay = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$zh.Hyperlinks.Item(2).TextToDisplay = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$zh.Hyperlinks.Item(3).TextToDisplay = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$zh.Hyperlinks.Item(4).TextToDisplay = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"

# ---------------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$de.Range("G2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.de-de.xlf"
$de.Range("H2").Value = "2016-09-06 17:44:27"
$de.Range("I2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$de.Range("J2").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.de-de.xlf"
$de.Range("K2").Value = "2016-09-06 17:44:48"

$de.Range("A3").Value = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$de.Range("G3").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.de-de.xlf"
$de.Range("H3").Value = "2016-09-06 17:44:27"
$de.Range("I3").Value = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$de.Range("J3").Value = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.de-de.xlf"
$de.Range("K3").Value = "2016-09-06 17:44:48"

$de.Hyperlinks.Item(1).TextToDisplay = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$de.Hyperlinks.Item(2).TextToDisplay = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$de.Hyperlinks.Item(3).TextToDisplay = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
$de.Hyperlinks.Item(4).TextToDisplay = "ffff6a949d48-3be0-4882-99e5-7ae793e0132c.md"
